# affichage dans backend corrigé
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: this reservation slot now belongs to "ahmed" ---
$ws.Cells.Item(4, 1).Value = "ahmed"
$ws.Cells.Item(4, 2).NumberFormat = "@"
$ws.Cells.Item(4, 2).Value = "2025-02-26"
$ws.Cells.Item(4, 3).Value = "En attente"

# --- Row 5: date corrected ---
$ws.Cells.Item(5, 2).NumberFormat = "@"
$ws.Cells.Item(5, 2).Value = "2025-02-25"

# --- Row 6: timestamp corrected ---
$ws.Cells.Item(6, 2).Value = "2025-03-06T12:34:49.917466400"

# --- New confirmed reservations for karoui (rows 7-13) ---
$newTimestamps = @(
    "2025-03-06T12:40:07.667345",
    "2025-03-06T12:44:21.003296400",
    "2025-03-06T12:44:34.228014700",
    "2025-03-06T12:53:51.661579400",
    "2025-03-06T12:58:30.918849200",
    "2025-03-06T13:01:30.447345700",
    "2025-03-06T13:02:13.757866900"
)

$r = 7
foreach ($ts in $newTimestamps) {
    # Clone row 6's formatting + values (A: "karoui", B: centered date, C: "Confirmé" fill) onto the new row
    $ws.Range("A6:C6").Copy($ws.Range("A" + $r + ":C" + $r))

    # Then patch in this row's own timestamp
    $ws.Cells.Item($r, 2).Value = $ts

    $r++
}
